# pl-2(完成版).pptx -- "Add files via upload" edit
#
# 1) Slide 3 speaker notes: remove the "Evaluate の部分は詳しい説明が必要" note
#    (it moves to slide 4's notes instead - see below).
# 2) Slide 4 speaker notes: replace the same placeholder note with the real
#    explanation about the SKET app vs. Yahoo! 知恵袋.
# 3) Slide 9, body placeholder: merge the two runs that make up
#    "（質問者は、解答者からの" + "解答を" into a single run, and keep the
#    red "評価" run immediately after it.

$p = $ppt.ActivePresentation

# --- 1) Slide 3 notes: drop the "Evaluate ..." reminder -------------------
$notes3 = $p.Slides.Item(3).NotesPage
$notes3.Shapes.Item(2).TextFrame.TextRange.Text = ""

# --- 2) Slide 4 notes: write the real explanation --------------------------
$notes4 = $p.Slides.Item(4).NotesPage
$notes4.Shapes.Item(2).TextFrame.TextRange.Text = "これだけだとyahoo知恵袋と同じだから以降で説明する特長で違いをあきらかにする。"

# --- 3) Slide 9: merge "（質問者は、解答者からの" + "解答を" into one run ----
$s9 = $p.Slides.Item(9)
$shp = $s9.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
# Characters(Start, Length) addresses the 12 + 3 = 15 characters that used to
# be split across two runs; re-assigning them as one string merges the runs
# (the surviving run keeps the formatting of the first of the two, which
# already matches the target rPr/dirty="0").
$merged = $tr.Characters(47, 15)
$merged.Text = "（質問者は、解答者からの解答を"
